$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 10871373
$ws.Range("I40").Value = 1784.1471
$ws.Range("K40").Value = 1784.1471
$ws.Range("M40").Value = -1609.1471
$ws.Range("H62").Value = 23812352
$ws.Range("I62").Value = 23812352
$ws.Range("K62").Value = 23812352
$ws.Range("M62").Value = -23811728
$ws.Range("H65").Value = 23812352
$ws.Range("I65").Value = 23812352
$ws.Range("K65").Value = 119061760
$ws.Range("M65").Value = -119058640
$ws.Range("H74").Value = 4700.5
$ws.Range("I74").Value = 4418
$ws.Range("J74").Value = 5045.778
$ws.Range("K74").Value = 4418
$ws.Range("L74").Value = 5045.778
$ws.Range("M74").Value = -3482
$ws.Range("N74").Value = -6917.778
$ws.Range("H77").Value = 4700.5
$ws.Range("I77").Value = 4418
$ws.Range("J77").Value = 5045.778
$ws.Range("K77").Value = 22090
$ws.Range("L77").Value = 25228.89
$ws.Range("M77").Value = -17410
$ws.Range("N77").Value = -34588.89
$ws.Range("H133").Value = 47772.5
$ws.Range("J133").Value = 47772.5
$ws.Range("L133").Value = 47772.5
$ws.Range("N133").Value = -57892.5
$ws.Range("H136").Value = 55485
$ws.Range("J136").Value = 55485
$ws.Range("L136").Value = 55485
$ws.Range("N136").Value = -65685
$ws.Range("H137").Value = 1294.25
$ws.Range("I137").Value = 1373.5217
$ws.Range("J137").Value = 1154
$ws.Range("K137").Value = 4120.5651
$ws.Range("L137").Value = 3462
$ws.Range("M137").Value = -1570.5651
$ws.Range("N137").Value = -8562
$ws.Range("H138").Value = 3481.6282
$ws.Range("I138").Value = 1682.5312
$ws.Range("J138").Value = 4733.174
$ws.Range("K138").Value = 5047.5936
$ws.Range("L138").Value = 14199.522
$ws.Range("M138").Value = 92.40639999999985
$ws.Range("N138").Value = -24479.522

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 134185.06
$ws.Range("I2").Value = 154559.77
$ws.Range("J2").Value = 1749.5
$ws.Range("K2").Value = 154559.77
$ws.Range("L2").Value = 1749.5
$ws.Range("M2").Value = -154446.77
$ws.Range("N2").Value = -1975.5
$ws.Range("H7").Value = 50712
$ws.Range("J7").Value = 50712
$ws.Range("L7").Value = 50712
$ws.Range("N7").Value = -50940
$ws.Range("H32").Value = 15220.719
$ws.Range("I32").Value = 11579.397
$ws.Range("K32").Value = 11579.397
$ws.Range("M32").Value = -11292.397
$ws.Range("H110").Value = 8104.16
$ws.Range("I110").Value = 8979.857
$ws.Range("K110").Value = 8979.857
$ws.Range("M110").Value = -6934.857
$ws.Range("H116").Value = 134185.06
$ws.Range("I116").Value = 154559.77
$ws.Range("J116").Value = 1749.5
$ws.Range("K116").Value = 154559.77
$ws.Range("L116").Value = 1749.5
$ws.Range("M116").Value = -152265.77
$ws.Range("N116").Value = -6337.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 134185.06
$ws.Range("I3").Value = 154559.77
$ws.Range("J3").Value = 1749.5
$ws.Range("K3").Value = 154559.77
$ws.Range("L3").Value = 1749.5
$ws.Range("M3").Value = -154445.77
$ws.Range("N3").Value = -1977.5
$ws.Range("H80").Value = 257.41666
$ws.Range("J80").Value = 333.41177
$ws.Range("L80").Value = 333.41177
$ws.Range("N80").Value = -2329.41177
$ws.Range("H83").Value = 257.41666
$ws.Range("J83").Value = 333.41177
$ws.Range("L83").Value = 1667.05885
$ws.Range("N83").Value = -11651.05885
$ws.Range("H133").Value = 44754.285
$ws.Range("J133").Value = 45546.668
$ws.Range("L133").Value = 45546.668
$ws.Range("N133").Value = -55666.668
$ws.Range("H134").Value = 23931.488
$ws.Range("I134").Value = 1666.9
$ws.Range("J134").Value = 202048.2
$ws.Range("K134").Value = 5000.700000000001
$ws.Range("L134").Value = 606144.6000000001
$ws.Range("M134").Value = -2465.700000000001
$ws.Range("N134").Value = -611214.6000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6761470
$ws.Range("I31").Value = 1397.8286
$ws.Range("J31").Value = 12828202
$ws.Range("K31").Value = 1397.8286
$ws.Range("L31").Value = 12828202
$ws.Range("M31").Value = -1102.8286
$ws.Range("N31").Value = -12828792
$ws.Range("H34").Value = 6761470
$ws.Range("I34").Value = 1397.8286
$ws.Range("J34").Value = 12828202
$ws.Range("K34").Value = 1397.8286
$ws.Range("L34").Value = 12828202
$ws.Range("M34").Value = -1195.8286
$ws.Range("N34").Value = -12828606
$ws.Range("H94").Value = 5145.8184
$ws.Range("J94").Value = 5108.6665
$ws.Range("L94").Value = 5108.6665
$ws.Range("N94").Value = -6010.6665
$ws.Range("H132").Value = 7409289.5
$ws.Range("I132").Value = 10001549
$ws.Range("J132").Value = 2834.5715
$ws.Range("K132").Value = 30004647
$ws.Range("L132").Value = 8503.7145
$ws.Range("M132").Value = -30002117
$ws.Range("N132").Value = -13563.7145

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 2274022.8
$ws.Range("I5").Value = 589.7931
$ws.Range("J5").Value = 8267618.5
$ws.Range("K5").Value = 1769.3793
$ws.Range("L5").Value = 24802855.5
$ws.Range("M5").Value = -1657.3793
$ws.Range("N5").Value = -24803079.5
$ws.Range("H121").Value = 744.7778
$ws.Range("I121").Value = 388.5
$ws.Range("J121").Value = 881.8077
$ws.Range("K121").Value = 1165.5
$ws.Range("L121").Value = 2645.4231
$ws.Range("M121").Value = 144.5
$ws.Range("N121").Value = -5265.4231
$ws.Range("H123").Value = 7451.4287
$ws.Range("I123").Value = 4386.6665
$ws.Range("K123").Value = 13159.9995
$ws.Range("M123").Value = -10709.9995
$ws.Range("H131").Value = 1819217.2
$ws.Range("J131").Value = 1162.425
$ws.Range("L131").Value = 3487.275
$ws.Range("N131").Value = -13567.275
$ws.Range("H132").Value = 4446775.5
$ws.Range("I132").Value = 1389.909
$ws.Range("J132").Value = 7939578.5
$ws.Range("K132").Value = 12509.181
$ws.Range("L132").Value = 71456206.5
$ws.Range("M132").Value = -9979.181
$ws.Range("N132").Value = -71461266.5
$ws.Range("H135").Value = 2274022.8
$ws.Range("I135").Value = 589.7931
$ws.Range("J135").Value = 8267618.5
$ws.Range("K135").Value = 5308.1379
$ws.Range("L135").Value = 74408566.5
$ws.Range("M135").Value = -2773.1379
$ws.Range("N135").Value = -74413636.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 58824500
$ws.Range("I113").Value = 111111890
$ws.Range("J113").Value = 1187.5
$ws.Range("K113").Value = 111111890
$ws.Range("L113").Value = 1187.5
$ws.Range("M113").Value = -111109720
$ws.Range("N113").Value = -5527.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 1454.375
$ws.Range("I61").Value = 1251.1538
$ws.Range("J61").Value = 2335
$ws.Range("K61").Value = 1251.1538
$ws.Range("L61").Value = 2335
$ws.Range("M61").Value = -1049.1538
$ws.Range("N61").Value = -2739
$ws.Range("H113").Value = 1454.375
$ws.Range("I113").Value = 1251.1538
$ws.Range("J113").Value = 2335
$ws.Range("K113").Value = 1251.1538
$ws.Range("L113").Value = 2335
$ws.Range("M113").Value = 918.8462
$ws.Range("N113").Value = -6675

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 2926045.8
$ws.Range("I136").Value = 1880.7142
$ws.Range("J136").Value = 11113708
$ws.Range("K136").Value = 5642.142599999999
$ws.Range("L136").Value = 33341124
$ws.Range("M136").Value = -3092.142599999999
$ws.Range("N136").Value = -33346224
